$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 1 header: B1/C1 now point at data1/data2 (string reuse, content unchanged) ---
$ws.Range("B1").Value = "data1"
$ws.Range("C1").Value = "data2"

# --- Row 3: customer name changed from "Deep" to "Deep74" ---
$ws.Range("D3").Value = "Deep74"

# --- Row 4: rename testcase id, rewrite opportunity fields ---
$ws.Range("A4").Value = "CreateOpportunity_ID"
$ws.Range("D4").Value = "500 Keyboard and computer"
$ws.Range("E4").Value = "100000"
$ws.Range("E4").NumberFormat = "@"

# --- Row 5: DeleteCustomer_ID ---
$ws.Range("A5").Value = "DeleteCustomer_ID"
$ws.Range("B5").Value = "pradeep@contentstudio.co.in"
$ws.Range("C5").Value = 8527529100
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:pradeep@contentstudio.co.in")
$ws.Range("B5").Style = $ws.Range("B2").Style

# --- Row 6: CreatingSalesTeam_ID ---
$ws.Range("A6").Value = "CreatingSalesTeam_ID"
$ws.Range("B6").Value = "pradeep@contentstudio.co.in"
$ws.Range("C6").Value = 8527529100
$ws.Range("D6").Value = "Sales965128"
$ws.Range("E6").Value = "pradeepsharma8850@gmail.com"
$ws.Range("F6").Value = "SalesTeam"
$ws.Range("G6").Value = "SalesCall"
$ws.Range("H6").Value = "SalesSMS"
$ws.Range("I6").Value = "SalesMeeting"
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:pradeep@contentstudio.co.in")
$ws.Range("B6").Style = $ws.Range("B2").Style
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:pradeepsharma8850@gmail.com")
$ws.Range("E6").Style = $ws.Range("B2").Style

# --- Row 7: CreateActivity_ID ---
$ws.Range("A7").Value = "CreateActivity_ID"
$ws.Range("B7").Value = "pradeep@contentstudio.co.in"
$ws.Range("C7").Value = 8527529100
$ws.Range("D7").Value = "AutomationMeeting"
$ws.Range("E7").Value = "pradeepsharma8850@gmail.com"
$ws.Range("F7").Value = "Automation Meeting with Customer"
$ws.Range("G7").Value = "Meeting"
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:pradeep@contentstudio.co.in")
$ws.Range("B7").Style = $ws.Range("B2").Style
$ws.Hyperlinks.Add($ws.Range("E7"), "mailto:pradeepsharma8850@gmail.com")
$ws.Range("E7").Style = $ws.Range("B2").Style

# --- Column widths ---
$ws.Columns.Item(4).ColumnWidth = 26.140625
$ws.Columns.Item(5).ColumnWidth = 30.5703125
$ws.Columns.Item(6).ColumnWidth = 33.5703125
$ws.Columns.Item(9).ColumnWidth = 13.140625

# --- View: scroll / selection ---
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("G9").Select()

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
